$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that would otherwise be auto-typed as numbers/dates by Excel's
# input parser (phone-like strings, JSHIR codes, passport-like numeric
# strings and the date column) must be pre-formatted as Text so they are
# stored verbatim, matching the source data ("Sana" column is a literal
# "yyyy-mm-dd" string here, not a date serial).
$textCols = @("F", "I", "J", "K")

$rows = @(
    @{
        A = "Nurmatova Rayxon Toxir qizi"
        B = "Mehnat muhofazasi va texnika xavfsizligi"
        C = "Rus tili"
        D = "Kunduzgi"
        E = "AD1289454"
        F = "62203066600043"
        G = "Toshkent shahri"
        H = "Mirzo Ulugʻbek tumani"
        I = "998947763216"
        J = "+998333906090"
        K = "2025-07-15"
    },
    @{
        A = "Yuldashev Doston Hamza o'g'li"
        B = "Mehnat muhofazasi va texnika xavfsizligi"
        C = "Rus tili"
        D = "Kunduzgi"
        E = "AB5180281"
        F = "52509006600014"
        G = "Toshkent shahri"
        H = "Mirzo Ulugʻbek tumani"
        I = "998997227826"
        J = "+998976010571"
        K = "2025-07-15"
    },
    @{
        A = "Ziyatov Ahror Ikrom ug'li"
        B = "Psixologiya"
        C = "O'zbek tili"
        D = "Kunduzgi"
        E = "AD8979548"
        F = "53006076230033"
        G = "Surxondaryo viloyati"
        H = "Denov tumani"
        I = "998200142003"
        J = "+998200142003"
        K = "2025-07-15"
    }
)

$startRow = 200
$cols = @("A","B","C","D","E","F","G","H","I","J","K")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Pre-format the numeric-looking / date columns as Text so Excel
    # stores the literal characters instead of coercing to a number or
    # a date serial (this mirrors pre-formatting the column as Text
    # before typing the values in).
    foreach ($c in $textCols) {
        $ws.Range("$c$r").NumberFormat = "@"
    }

    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $data[$c]
    }

    # Drop the temporary Text number format again so the new rows keep
    # the same (default / no explicit style) look as the rest of the
    # sheet's data rows.
    foreach ($c in $textCols) {
        $ws.Range("$c$r").ClearFormats()
    }
}

